$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 250
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 250
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H98").Value = 4392.136
$ws.Range("I98").Value = 3101.8
$ws.Range("J98").Value = 7157.143
$ws.Range("K98").Value = 3101.8
$ws.Range("L98").Value = 7157.143
$ws.Range("M98").Value = -1603.8
$ws.Range("N98").Value = -10153.143

$ws.Range("H122").Value = 4392.136
$ws.Range("I122").Value = 3101.8
$ws.Range("J122").Value = 7157.143
$ws.Range("K122").Value = 9305.400000000001
$ws.Range("L122").Value = 21471.429
$ws.Range("M122").Value = -6855.400000000001
$ws.Range("N122").Value = -26371.429

$ws.Range("H131").Value = 2748.125
$ws.Range("I131").Value = 2340.7144
$ws.Range("J131").Value = 5600
$ws.Range("K131").Value = 7022.1432
$ws.Range("L131").Value = 16800
$ws.Range("M131").Value = -1982.1432
$ws.Range("N131").Value = -26880

$ws.Range("H132").Value = 23907226
$ws.Range("I132").Value = 31377034
$ws.Range("J132").Value = 3844
$ws.Range("K132").Value = 94131102
$ws.Range("L132").Value = 11532
$ws.Range("M132").Value = -94128572
$ws.Range("N132").Value = -16592

$ws.Range("H138").Value = 2140.09
$ws.Range("I138").Value = 1057.7435
$ws.Range("J138").Value = 2832.082
$ws.Range("K138").Value = 3173.2305
$ws.Range("L138").Value = 8496.245999999999
$ws.Range("M138").Value = 1966.7695
$ws.Range("N138").Value = -18776.246

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1666.3334
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -884

$ws.Range("H6").Value = 9869.75
$ws.Range("J6").Value = 10497
$ws.Range("L6").Value = 10497
$ws.Range("N6").Value = -10843

$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 10000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 10000
$ws.Range("N14").Value = 0
$ws.Range("M14").Value = -9825
$ws.Range("L14").ClearContents()

$ws.Range("H32").Value = 12080.275
$ws.Range("I32").Value = 8475.146000000001
$ws.Range("K32").Value = 8475.146000000001
$ws.Range("M32").Value = -8188.146000000001

$ws.Range("I102").Value = 2300
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2300
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = -678
$ws.Range("M102").ClearContents()

$ws.Range("H137").Value = 53570
$ws.Range("J137").Value = 53570
$ws.Range("L137").Value = 53570
$ws.Range("N137").Value = -63770

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 606.75
$ws.Range("I25").Value = 606.75
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 606.75
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = -371.75
$ws.Range("M25").ClearContents()

$ws.Range("H134").Value = 3123.262
$ws.Range("I134").Value = 1714.5
$ws.Range("J134").Value = 7631.3
$ws.Range("K134").Value = 5143.5
$ws.Range("L134").Value = 22893.9
$ws.Range("M134").Value = -2608.5
$ws.Range("N134").Value = -27963.9

$ws.Range("H137").Value = 35361.25
$ws.Range("J137").Value = 35361.25
$ws.Range("L137").Value = 35361.25
$ws.Range("N137").Value = -45561.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2073.4849
$ws.Range("I58").Value = 1776.8276
$ws.Range("J58").Value = 4224.25
$ws.Range("K58").Value = 1776.8276
$ws.Range("L58").Value = 4224.25
$ws.Range("M58").Value = -1573.8276
$ws.Range("N58").Value = -4630.25

$ws.Range("H136").Value = 2073.4849
$ws.Range("I136").Value = 1776.8276
$ws.Range("J136").Value = 4224.25
$ws.Range("K136").Value = 5330.4828
$ws.Range("L136").Value = 12672.75
$ws.Range("M136").Value = -2780.4828
$ws.Range("N136").Value = -17772.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 253251.25
$ws.Range("I55").Value = 1000000
$ws.Range("J55").Value = 4335
$ws.Range("K55").Value = 3000000
$ws.Range("L55").Value = 13005
$ws.Range("M55").Value = -2999823
$ws.Range("N55").Value = -13359

$ws.Range("H131").Value = 12821500
$ws.Range("J131").Value = 1043.2941
$ws.Range("L131").Value = 3129.8823
$ws.Range("N131").Value = -13209.8823

$ws.Range("H132").Value = 2735.88
$ws.Range("I132").Value = 1114.7693
$ws.Range("J132").Value = 4492.0835
$ws.Range("K132").Value = 10032.9237
$ws.Range("L132").Value = 40428.7515
$ws.Range("M132").Value = -7502.923699999999
$ws.Range("N132").Value = -45488.7515

$ws.Range("H134").Value = 2987.8206
$ws.Range("I134").Value = 1811.875
$ws.Range("J134").Value = 4869.3335
$ws.Range("K134").Value = 5435.625
$ws.Range("L134").Value = 14608.0005
$ws.Range("M134").Value = -365.625
$ws.Range("N134").Value = -24748.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25002558
$ws.Range("J80").Value = 2939.8572
$ws.Range("L80").Value = 2939.8572
$ws.Range("N80").Value = -4935.8572

$ws.Range("H83").Value = 25002558
$ws.Range("J83").Value = 2939.8572
$ws.Range("L83").Value = 14699.286
$ws.Range("N83").Value = -24683.286

$ws.Range("H126").Value = 3834.8062
$ws.Range("I126").Value = 2810.75
$ws.Range("J126").Value = 5200.2144
$ws.Range("K126").Value = 8432.25
$ws.Range("L126").Value = 15600.6432
$ws.Range("M126").Value = -5962.25
$ws.Range("N126").Value = -20540.6432

$ws.Range("H132").Value = 3201.476
$ws.Range("I132").Value = 1326.3
$ws.Range("K132").Value = 3978.9
$ws.Range("M132").Value = -1448.9

$ws.Range("H137").Value = 84827.5
$ws.Range("J137").Value = 84827.5
$ws.Range("L137").Value = 84827.5
$ws.Range("N137").Value = -95027.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 499.4138
$ws.Range("I16").Value = 507.92593
$ws.Range("K16").Value = 507.92593
$ws.Range("M16").Value = -337.92593

$ws.Range("H133").Value = 53929.91
$ws.Range("J133").Value = 53929.91
$ws.Range("L133").Value = 53929.91
$ws.Range("N133").Value = -58989.91

$ws.Range("H136").Value = 3612.5757
$ws.Range("I136").Value = 1847.6471
$ws.Range("J136").Value = 5487.8125
$ws.Range("K136").Value = 5542.9413
$ws.Range("L136").Value = 16463.4375
$ws.Range("M136").Value = -2992.9413
$ws.Range("N136").Value = -21563.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18908556
$ws.Range("I81").Value = 18908556
$ws.Range("K81").Value = 37817112
$ws.Range("M81").Value = -37816051

$ws.Range("H84").Value = 18908556
$ws.Range("I84").Value = 18908556
$ws.Range("K84").Value = 189085560
$ws.Range("M84").Value = -189080256

$ws.Range("H93").Value = 34800
$ws.Range("J93").Value = 34800
$ws.Range("L93").Value = 34800
$ws.Range("N93").Value = -39792

$ws.Range("H123").Value = 29886
$ws.Range("J123").Value = 29886
$ws.Range("L123").Value = 29886
$ws.Range("N123").Value = -39686
